$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.307.75'
$ws.Range('E2').Value = '  -2.93%  '
$ws.Range('D3').Value = '3.426.14'
$ws.Range('E3').Value = '  -5.35%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.44'
$ws.Range('E5').Value = '  -5.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '188.18'
$ws.Range('E6').Value = '  -4.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.605'
$ws.Range('E7').Value = '  -3.35%  '
$ws.Range('D8').Value = '3.411.95'
$ws.Range('E8').Value = '  -5.40%  '
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.199'
$ws.Range('E10').Value = '  -6.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.609'
$ws.Range('E11').Value = '  -5.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '50.43'
$ws.Range('E12').Value = '  -5.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000280'
$ws.Range('E13').Value = '  -7.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.97'
$ws.Range('E14').Value = '  -6.02%  '
$ws.Range('D15').Value = '3.980.80'
$ws.Range('E15').Value = '  -5.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '628.63'
$ws.Range('E16').Value = '  +3.75%  '
$ws.Range('D17').Value = '68.331.44'
$ws.Range('E17').Value = '  -3.02%  '
$ws.Range('D18').Value = '3.436.49'
$ws.Range('E18').Value = '  -5.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.06'
$ws.Range('E20').Value = '  -7.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.93'
$ws.Range('E21').Value = '  -5.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.930'
$ws.Range('E22').Value = '  -6.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.58'
$ws.Range('E23').Value = '  -2.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.29'
$ws.Range('E24').Value = '  +1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.85'
$ws.Range('E25').Value = '  -5.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.21'
$ws.Range('E26').Value = '  -8.61%  '
$ws.Range('E27').Value = '  +1.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.80'
$ws.Range('E28').Value = '  -5.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.71'
$ws.Range('E29').Value = '  -8.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.08'
$ws.Range('E30').Value = '  -6.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.93'
$ws.Range('E31').Value = '  -5.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.09'
$ws.Range('E32').Value = '  -12.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.59'
$ws.Range('E33').Value = '  -9.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.43'
$ws.Range('E34').Value = '  -6.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '60.25'
$ws.Range('E35').Value = '  -4.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.107'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').Value = '3.615.47'
$ws.Range('E38').Value = '  -8.18%  '
$ws.Range('D39').Value = '0.0₃0771'
$ws.Range('E39').Value = '  -13.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '494.55'
$ws.Range('E40').Value = '  -5.34%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.43'
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.83'
$ws.Range('E42').Value = '  -7.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.363'
$ws.Range('E43').Value = '  -6.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.131'
$ws.Range('E44').Value = '  -3.57%  '
$ws.Range('B45').Value = 'CoreDAO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.45'
$ws.Range('E45').Value = '  +67.21%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '33.74'
$ws.Range('E46').Value = '  -7.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0432'
$ws.Range('E47').Value = '  -6.33%  '
$ws.Range('E48').Value = '  -5.49%  '
$ws.Range('E49').Value = '  -4.24%  '
$ws.Range('E50').Value = '  -5.09%  '
$ws.Range('E51').Value = '  -0.40%  '
